# Horarios actualizados Linea 141 - scrape refresh 19:35:56
# Updates the "LP1912" schedule sheet with refreshed arrival-time scrape
# data (header counters + reordered/updated rows + 7 newly scraped rows
# appended at the end of the table), and refreshes the "last updated"
# timestamp on the two companion sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# ---- Sheet: LP1912 ----
$ws1.Range("A2").Value = "Última actualización: 19:35:56"
$ws1.Range("A3").Value = "Total filas: 340"
$ws1.Range("A56").Value = "07:38:39"
$ws1.Range("C56").Value = "14_ABASTO"
$ws1.Range("D56").Value = 99
$ws1.Range("A57").Value = "08:27:16"
$ws1.Range("C57").Value = "27_EL RETIRO"
$ws1.Range("D57").Value = 50
$ws1.Range("A111").Value = "11:52:01"
$ws1.Range("C111").Value = "225_GOMEZ"
$ws1.Range("D111").Value = 2
$ws1.Range("A112").Value = "11:54:18"
$ws1.Range("C112").Value = "15X38_ABASTO"
$ws1.Range("D112").Value = 0
$ws1.Range("A113").Value = "10:50:41"
$ws1.Range("C113").Value = "23_HERNANDEZ"
$ws1.Range("D113").Value = 64
$ws1.Range("C120").Value = "15_ABASTO"
$ws1.Range("C121").Value = "16_P MOR-SANTA ANA"
$ws1.Range("A142").Value = "11:11:33"
$ws1.Range("C142").Value = "15X38_ABASTO"
$ws1.Range("D142").Value = 97
$ws1.Range("A143").Value = "10:50:41"
$ws1.Range("C143").Value = "16_SANTA ANA"
$ws1.Range("D143").Value = 118
$ws1.Range("A144").Value = "11:47:17"
$ws1.Range("C144").Value = "14_ABASTO"
$ws1.Range("D144").Value = 61
$ws1.Range("A216").Value = "14:12:26"
$ws1.Range("C216").Value = "14_ABASTO"
$ws1.Range("D216").Value = 113
$ws1.Range("A217").Value = "15:17:33"
$ws1.Range("C217").Value = "16_SANTA ANA"
$ws1.Range("D217").Value = 48
$ws1.Range("A257").Value = "16:52:42"
$ws1.Range("C257").Value = "23_HERNANDEZ"
$ws1.Range("D257").Value = 54
$ws1.Range("A258").Value = "15:58:05"
$ws1.Range("C258").Value = "215_EL PELIGRO"
$ws1.Range("D258").Value = 108
$ws1.Range("C259").Value = "215_EL PELIGRO"
$ws1.Range("C261").Value = "215B_EL PATO"
$ws1.Range("A302").Value = "18:12:30"
$ws1.Range("C302").Value = "27_EL RETIRO"
$ws1.Range("D302").Value = 64
$ws1.Range("A304").Value = "18:44:57"
$ws1.Range("C304").Value = "14X44_ABASTO"
$ws1.Range("D304").Value = 32
$ws1.Range("A323").Value = "19:35:56"
$ws1.Range("B323").Value = "20:04"
$ws1.Range("C323").Value = "23_HERNANDEZ"
$ws1.Range("D323").Value = 29
$ws1.Range("A324").Value = "18:44:57"
$ws1.Range("B324").Value = "20:06"
$ws1.Range("D324").Value = 82
$ws1.Range("B325").Value = "20:07"
$ws1.Range("C325").Value = "215C_EL PATO"
$ws1.Range("D325").Value = 75
$ws1.Range("A326").Value = "18:52:02"
$ws1.Range("B326").Value = "20:08"
$ws1.Range("D326").Value = 76
$ws1.Range("B327").Value = "20:09"
$ws1.Range("C327").Value = "23_HERNANDEZ"
$ws1.Range("D327").Value = 85
$ws1.Range("A328").Value = "18:44:57"
$ws1.Range("C328").Value = "14_ABASTO"
$ws1.Range("D328").Value = 88
$ws1.Range("A329").Value = "18:31:25"
$ws1.Range("B329").Value = "20:12"
$ws1.Range("C329").Value = "215C_EL PATO"
$ws1.Range("D329").Value = 101
$ws1.Range("A330").Value = "18:52:02"
$ws1.Range("B330").Value = "20:13"
$ws1.Range("C330").Value = "14_ABASTO"
$ws1.Range("D330").Value = 81
$ws1.Range("A331").Value = "18:44:57"
$ws1.Range("B331").Value = "20:21"
$ws1.Range("D331").Value = 97
$ws1.Range("A332").Value = "18:31:25"
$ws1.Range("B332").Value = "20:22"
$ws1.Range("C332").Value = "15_ABASTO"
$ws1.Range("D332").Value = 111
$ws1.Range("A333").Value = "18:44:57"
$ws1.Range("B333").Value = "20:30"
$ws1.Range("D333").Value = 106
$ws1.Range("A334").Value = "18:52:02"
$ws1.Range("B334").Value = "20:31"
$ws1.Range("C334").Value = "10_OLMOS"
$ws1.Range("D334").Value = 99
$ws1.Range("A335").Value = "19:35:56"
$ws1.Range("B335").Value = "20:33"
$ws1.Range("C335").Value = "16_SANTA ANA"
$ws1.Range("D335").Value = 58
$ws1.Range("A336").Value = "19:35:56"
$ws1.Range("B336").Value = "20:42"
$ws1.Range("C336").Value = "17_ROMERO"
$ws1.Range("D336").Value = 67
$ws1.Range("B337").Value = "20:43"
$ws1.Range("C337").Value = "17_ROMERO"
$ws1.Range("D337").Value = 92
$ws1.Range("B338").Value = "20:47"
$ws1.Range("C338").Value = "215B_EL PATO"
$ws1.Range("D338").Value = 96
$ws1.Range("A339").Value = "18:52:02"
$ws1.Range("B339").Value = "20:48"
$ws1.Range("C339").Value = "215B_EL PATO"
$ws1.Range("D339").Value = 116
$ws1.Range("E339").Value = "LP1912"
$ws1.Range("A340").Value = "19:11:45"
$ws1.Range("B340").Value = "20:56"
$ws1.Range("C340").Value = "27_EL RETIRO"
$ws1.Range("D340").Value = 105
$ws1.Range("E340").Value = "LP1912"
$ws1.Range("A341").Value = "19:35:56"
$ws1.Range("B341").Value = "20:57"
$ws1.Range("C341").Value = "23_HERNANDEZ"
$ws1.Range("D341").Value = 82
$ws1.Range("E341").Value = "LP1912"
$ws1.Range("A342").Value = "19:11:45"
$ws1.Range("B342").Value = "21:06"
$ws1.Range("C342").Value = "10_OLMOS"
$ws1.Range("D342").Value = 115
$ws1.Range("E342").Value = "LP1912"
$ws1.Range("A343").Value = "19:35:56"
$ws1.Range("B343").Value = "21:09"
$ws1.Range("C343").Value = "15_ABASTO"
$ws1.Range("D343").Value = 94
$ws1.Range("E343").Value = "LP1912"
$ws1.Range("A344").Value = "19:35:56"
$ws1.Range("B344").Value = "21:28"
$ws1.Range("C344").Value = "11_ETCHEVERRY"
$ws1.Range("D344").Value = 113
$ws1.Range("E344").Value = "LP1912"
$ws1.Range("A345").Value = "19:35:56"
$ws1.Range("B345").Value = "21:33"
$ws1.Range("C345").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D345").Value = 118
$ws1.Range("E345").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws2.Range("A2").Value = "Última actualización: 19:35:56"

# ---- Sheet: 6203-6173 ----
$ws3.Range("A2").Value = "Última actualización: 19:35:56"
